$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("McCauley rotifers")

# --- Restructure "McCauley rotifers": insert a new column A for species,
# pushing the existing genus/constant/formula columns from A,B,C to B,C,D ---
$ws.Columns("A").Insert()

# Header row
$ws.Range("A1").Value = "species"
$ws.Range("B1").Value = "genus"
$ws.Range("C1").Value = "constant"
$ws.Range("D1").Value = "formula"

# Capitalise every genus name now sitting in column B (rows 2-21)
for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Range("B$r")
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cap = $val.Substring(0,1).ToUpper() + $val.Substring(1)
        $cell.Value = $cap
    }
}

# Rows 13 & 14 hold species-level names (two words) - move them from the
# genus column (B) into the new species column (A) and clear out B
$ws.Range("A13").Value = $ws.Range("B13").Value2
$ws.Range("B13").ClearContents()

$ws.Range("A14").Value = $ws.Range("B14").Value2
$ws.Range("B14").ClearContents()

# --- View/selection updates on this sheet ---
$ws.Activate()
$ws.Application.ActiveWindow.Zoom = 140
$ws.Range("D14").Select()
